$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preserving its "General"
# number format and default (unstyled) appearance. Values in the source
# sheet are plain text (price/volume strings, some with a trailing "%"),
# so a direct .Value assignment would get auto-coerced into a number by
# Excel. Temporarily forcing a Text format before the assignment keeps the
# literal string, then resetting the style back to Normal removes the
# Text-format marker again so the cell ends up indistinguishable from the
# original (no number format / no explicit style).
function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "303.18"
Set-TextValue $ws "E2" "1.35%"
Set-TextValue $ws "D3" "32.72"
Set-TextValue $ws "E3" "4.55%"
Set-TextValue $ws "D4" "4.928"
Set-TextValue $ws "E4" "-3.40%"
Set-TextValue $ws "D5" "0.07835"
Set-TextValue $ws "E5" "-1.36%"
Set-TextValue $ws "D6" "2.033"
Set-TextValue $ws "E6" "-10.87%"
Set-TextValue $ws "D7" "7.835"
Set-TextValue $ws "E7" "0.77%"
Set-TextValue $ws "E8" "-1.50%"
Set-TextValue $ws "D9" "0.9228"
Set-TextValue $ws "E9" "-0.17%"
Set-TextValue $ws "E10" "1.17%"
Set-TextValue $ws "D11" "0.07877"
Set-TextValue $ws "E11" "5.11%"
Set-TextValue $ws "D12" "0.08641"
Set-TextValue $ws "E12" "-8.92%"
Set-TextValue $ws "D13" "0.03140"
Set-TextValue $ws "E13" "3.16%"
Set-TextValue $ws "E14" "0.17%"
Set-TextValue $ws "D15" "0.001509"
Set-TextValue $ws "E15" "-0.72%"
Set-TextValue $ws "D16" "0.005923"
Set-TextValue $ws "E16" "0.34%"
Set-TextValue $ws "D17" "3.465"
Set-TextValue $ws "E17" "-0.54%"
Set-TextValue $ws "E18" "-4.85%"
Set-TextValue $ws "E19" "1.10%"
Set-TextValue $ws "E20" "-3.34%"
Set-TextValue $ws "D21" "4.312"
Set-TextValue $ws "E21" "10.30%"
Set-TextValue $ws "D22" "0.1990"
Set-TextValue $ws "E22" "17.11%"
Set-TextValue $ws "D23" "0.04575"
Set-TextValue $ws "E23" "-0.84%"
Set-TextValue $ws "E24" "-2.01%"
Set-TextValue $ws "D25" "0.004449"
Set-TextValue $ws "E25" "-0.77%"
Set-TextValue $ws "E26" "4.22%"
Set-TextValue $ws "D39" "0.01739"
Set-TextValue $ws "E39" "-1.02%"
Set-TextValue $ws "D40" "0.04781"
Set-TextValue $ws "E40" "3.54%"
Set-TextValue $ws "D41" "0.007478"
Set-TextValue $ws "E41" "7.23%"
Set-TextValue $ws "E42" "-0.28%"
Set-TextValue $ws "D43" "0.002339"
Set-TextValue $ws "E43" "6.91%"
Set-TextValue $ws "D44" "0.01054"
Set-TextValue $ws "E44" "2.75%"
Set-TextValue $ws "D45" "0.00006235"
Set-TextValue $ws "E45" "-2.06%"
Set-TextValue $ws "E46" "0.05%"
Set-TextValue $ws "E47" "-61.11%"
Set-TextValue $ws "D48" "0.8234"
Set-TextValue $ws "E48" "10.25%"
Set-TextValue $ws "D49" "0.00002099"
Set-TextValue $ws "E49" "0.05%"
Set-TextValue $ws "D50" "0.0001999"
Set-TextValue $ws "E50" "0.05%"
